$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was old row 4: Segunda, 44881)
$ws.Range("D2").Value = 44881
$ws.Range("L2").Value = "Segunda"
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 11250
$ws.Range("O2").Value = 11250
$ws.Range("P2").Value = 11250
$ws.Range("S2").Value = 11250

# Row 3 (was old row 5: Primera, 44874, vol 200)
$ws.Range("D3").Value = 44874
$ws.Range("M3").Value = 200
$ws.Range("P3").Value = 7750
$ws.Range("S3").Value = 7750

# Row 4 (was old row 2: Primera, 44923)
$ws.Range("D4").Value = 44923
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 80
$ws.Range("N4").Value = 7500
$ws.Range("O4").Value = 8000
$ws.Range("P4").Value = 7625
$ws.Range("S4").Value = 7625

# Row 5 (was old row 3: Primera, 44923, vol 80)
$ws.Range("D5").Value = 44923
$ws.Range("M5").Value = 80
$ws.Range("P5").Value = 7625
$ws.Range("S5").Value = 7625
